$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet shrinks from 25 to 23 data rows: drop the old trailing rows 24-25
# (their content -- "Requisitos:" / the LOT2056 note -- moves up into rows 22-23).
$ws.Rows("24:25").Delete()

# Row 13: "01/01/2020" must stay a literal text string (like the existing B8/C8 cells),
# not get reinterpreted as a date serial, so copy it instead of assigning the literal.
$ws.Range("B8:C8").Copy($ws.Range("B13"))

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Rows(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Development of the course conclusion work under orientation of a leader professor, which must deal with subjects related to Biochemical Engineering."
$ws.Range("C14").Value = "Development of the course conclusion work under orientation of a leader professor, which must deal with subjects related to Biochemical Engineering."
$ws.Rows(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C15").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Rows(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Elaboration of a monograph of course conclusion presenting: (1) the subject and its importance, (2) the objectives, (3) the bibliographic revision, (4) the scientific methodology, (5) the development of the project, (6) the analysis and discussion of the results, (7) the conclusion and recommendations for the future works and (8) references. The document must attend to the Brazilian Association of Technical Norms."
$ws.Range("C16").Value = "Elaboration of a monograph of course conclusion presenting: (1) the subject and its importance, (2) the objectives, (3) the bibliographic revision, (4) the scientific methodology, (5) the development of the project, (6) the analysis and discussion of the results, (7) the conclusion and recommendations for the future works and (8) references. The document must attend to the Brazilian Association of Technical Norms."
$ws.Rows(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows(17).AutoFit()

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8853480 - Tatiane da Franca Silva"
$ws.Range("C18").Value = "8853480 - Tatiane da Franca Silva"
$ws.Rows(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica."
$ws.Range("C19").Value = "Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica."
$ws.Rows(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota da disciplina será decidida pelos docentes da banca"
$ws.Range("C20").Value = "A nota da disciplina será decidida pelos docentes da banca"
$ws.Rows(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Reapresentação do trabalho modificado para nova avaliação."
$ws.Range("C21").Value = "Reapresentação do trabalho modificado para nova avaliação."
$ws.Rows(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows(22).AutoFit()

# Row 23
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOT2056 -  Trabalho de Conclusão de Curso I  (Requisito)`n"
$ws.Range("C23").Value = "LOT2056 -  Trabalho de Conclusão de Curso I  (Requisito)`n"
$ws.Rows(23).RowHeight = 30
